# Updated cryptos list (prices + 1h volume deltas), and swapped the
# Toncoin/Cosmos row order (rows 28-29), per the Dec 11 2023 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal TEXT (matches the source sheet, where
# every data cell - including bare-looking numbers like "242.39" - is stored
# as a string, not a numeric cell). Without forcing the "@" text format first,
# Excel auto-converts digit-and-dot strings to floating point numbers, which
# both changes the cell type and can introduce binary rounding (e.g. 242.39
# -> 242.38999999999999). Resetting the Style back to "Normal" afterwards drops
# the now-unneeded explicit number format so the cell keeps the workbook's
# original (unstyled) look.
function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "41.219.28"
$ws.Range("E2").Value = "  -6.02%  "

$ws.Range("D3").Value = "2.213.34"

$ws.Range("E4").Value = "  +0.11%  "

Set-TextValue "D5" "242.39"
$ws.Range("E5").Value = "  +1.05%  "

$ws.Range("E6").Value = "  -6.74%  "

Set-TextValue "D7" "69.88"
$ws.Range("E7").Value = "  -5.90%  "

$ws.Range("E8").Value = "  +0.21%  "

Set-TextValue "D9" "0.553"
$ws.Range("E9").Value = "  -7.63%  "

Set-TextValue "D10" "38.08"
$ws.Range("E10").Value = "  +2.80%  "

$ws.Range("E11").Value = "  -7.04%  "

Set-TextValue "D12" "57.87"
$ws.Range("E12").Value = "  -3.38%  "

$ws.Range("E13").Value = "  -3.45%  "

Set-TextValue "D14" "6.71"
$ws.Range("E14").Value = "  -7.95%  "

$ws.Range("D15").Value = "2.544.00"
$ws.Range("E15").Value = "  -6.39%  "

Set-TextValue "D16" "14.77"
$ws.Range("E16").Value = "  -9.69%  "

Set-TextValue "D17" "0.838"
$ws.Range("E17").Value = "  -9.77%  "

$ws.Range("D18").Value = "2.215.11"
$ws.Range("E18").Value = "  -6.62%  "

$ws.Range("D19").Value = "41.192.33"
$ws.Range("E19").Value = "  -6.00%  "

$ws.Range("D20").Value = "0.0₃0947"
$ws.Range("E20").Value = "  -8.31%  "

Set-TextValue "D21" "72.19"
$ws.Range("E21").Value = "  -6.87%  "

Set-TextValue "D22" "6.09"
$ws.Range("E22").Value = "  -8.01%  "

Set-TextValue "D23" "231.56"
$ws.Range("E23").Value = "  -8.78%  "

$ws.Range("E24").Value = "  +6.55%  "

$ws.Range("E25").Value = "  +0.14%  "

$ws.Range("E26").Value = "  -5.23%  "

Set-TextValue "D27" "2.41"
$ws.Range("E27").Value = "  -3.32%  "

$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D28" "9.76"
$ws.Range("E28").Value = "  -7.67%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D29" "2.18"
$ws.Range("E29").Value = "  -4.90%  "

Set-TextValue "D30" "172.32"
$ws.Range("E30").Value = "  -1.81%  "

Set-TextValue "D31" "20.44"
$ws.Range("E31").Value = "  -8.66%  "

$ws.Range("E32").Value = "  -8.19%  "

Set-TextValue "D33" "0.123"
$ws.Range("E33").Value = "  -8.44%  "

Set-TextValue "D34" "0.0708"
$ws.Range("E34").Value = "  -6.75%  "

Set-TextValue "D35" "5.20"
$ws.Range("E35").Value = "  -4.07%  "

$ws.Range("E36").Value = "  -9.75%  "

Set-TextValue "D37" "3.90"
$ws.Range("E37").Value = "  +2.40%  "

Set-TextValue "D38" "23.72"
$ws.Range("E38").Value = "  +15.41%  "

Set-TextValue "D39" "0.0278"
$ws.Range("E39").Value = "  -1.29%  "

$ws.Range("E40").Value = "  -5.42%  "

Set-TextValue "D41" "5.83"
$ws.Range("E41").Value = "  -11.92%  "

Set-TextValue "D42" "64.39"
$ws.Range("E42").Value = "  -1.73%  "

$ws.Range("E43").Value = "  -11.06%  "

$ws.Range("E44").Value = "  -3.26%  "

Set-TextValue "D45" "8.62"
$ws.Range("E45").Value = "  -5.05%  "

$ws.Range("E46").Value = "  -6.76%  "

$ws.Range("E47").Value = "  -0.08%  "

$ws.Range("E48").Value = "  +10.56%  "

Set-TextValue "D49" "4.47"
$ws.Range("E49").Value = "  +1.83%  "

$ws.Range("E50").Value = "  -6.34%  "

$ws.Range("E51").Value = "  -5.38%  "
